$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.007742333333333334
$ws.Range("N2").Value = 0.023227
$ws.Range("O2").Value = 0.001217676423630818
$ws.Range("P2").Value = 0.001217676423630818
$ws.Range("Q2").Value = 0.001701264195777778
$ws.Range("R2").Value = 0.015311377762
$ws.Range("S2").Value = 0.001217676423630818
$ws.Range("T2").Value = 0.001217676423630818

# Row 3
$ws.Range("O3").Value = 0.9497929577862038
$ws.Range("P3").Value = 0.9497929577862039
$ws.Range("S3").Value = 0.9497929577862038
$ws.Range("T3").Value = 0.9497929577862039

# Row 4
$ws.Range("M4").Value = 0.3114883333333334
$ws.Range("N4").Value = 0.9344650000000001
$ws.Range("O4").Value = 0.04898936579016542
$ws.Range("P4").Value = 0.04898936579016543
$ws.Range("Q4").Value = 0.06844499275444445
$ws.Range("R4").Value = 0.6160049347900001
$ws.Range("S4").Value = 0.04898936579016542
$ws.Range("T4").Value = 0.04898936579016543
